$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "mNGplate11_sorted_A10_PATZ1-N"
$ws.Range("D3").Value = "mNGplate11_sorted_A11_KDELR3-C"
$ws.Range("D4").Value = "mNGplate11_sorted_A12_MYH9-C"
